$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.67"
$ws.Range("E2").Value = "'-1.12%"
$ws.Range("D3").Value = "'35.74"
$ws.Range("E3").Value = "'-0.24%"
$ws.Range("D4").Value = "'5.033"
$ws.Range("E4").Value = "'-0.09%"
$ws.Range("D5").Value = "'0.07973"
$ws.Range("D6").Value = "'1.850"
$ws.Range("E6").Value = "'-5.74%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.768"
$ws.Range("E7").Value = "'-0.35%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9215"
$ws.Range("E8").Value = "'-1.17%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1276"
$ws.Range("E9").Value = "'-4.62%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1884"
$ws.Range("E10").Value = "'-2.17%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08955"
$ws.Range("E11").Value = "'-3.04%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03415"
$ws.Range("E12").Value = "'-2.66%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09840"
$ws.Range("E13").Value = "'-0.35%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001405"
$ws.Range("E14").Value = "'-2.03%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006247"
$ws.Range("E15").Value = "'7.24%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.850"
$ws.Range("E16").Value = "'6.89%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.123"
$ws.Range("E17").Value = "'-0.17%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'3.296"
$ws.Range("E18").Value = "'12.04%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3405"
$ws.Range("E19").Value = "'-0.69%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1340"
$ws.Range("E20").Value = "'0.64%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'4.803"
$ws.Range("E21").Value = "'-7.31%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2344"
$ws.Range("E22").Value = "'-9.64%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04384"
$ws.Range("E23").Value = "'-0.12%"
$ws.Range("D24").Value = "'0.001235"
$ws.Range("E24").Value = "'1.33%"
$ws.Range("D25").Value = "'0.004847"
$ws.Range("E25").Value = "'1.51%"
$ws.Range("D27").Value = "'0.0001302"
$ws.Range("E27").Value = "'-21.09%"
$ws.Range("E28").Value = "'42.45%"
$ws.Range("D39").Value = "'0.01920"
$ws.Range("E39").Value = "'-3.41%"
$ws.Range("D40").Value = "'0.05126"
$ws.Range("E40").Value = "'2.15%"
$ws.Range("D41").Value = "'0.007546"
$ws.Range("E41").Value = "'-1.07%"
$ws.Range("D42").Value = "'0.01016"
$ws.Range("E42").Value = "'-11.00%"
$ws.Range("D43").Value = "'0.1344"
$ws.Range("E43").Value = "'-2.51%"
$ws.Range("D44").Value = "'0.002114"
$ws.Range("E44").Value = "'0.87%"
$ws.Range("D45").Value = "'0.009879"
$ws.Range("E45").Value = "'-12.77%"
$ws.Range("D46").Value = "'0.00006192"
$ws.Range("E46").Value = "'-3.01%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.38%"
$ws.Range("D48").Value = "'64.85"
$ws.Range("E48").Value = "'-0.17%"
$ws.Range("D49").Value = "'0.001253"
$ws.Range("E49").Value = "'5.42%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.38%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.38%"
